$wb = $excel.ActiveWorkbook

# --- Fix the misspelled tournament name "Rolland Garros" -> "Roland Garros" ---
# Rename every occurrence across all three sheets (all share the same
# mis-spelling, which lived in the shared-strings table). Renaming *every*
# use lets the old shared string drop out and the corrected text gets
# appended as a fresh shared-string entry, matching the authoring tool's
# behaviour of moving the corrected entry to the end of the table.
foreach ($sheetName in @("Roger Federer", "Rafael Nadal", "Novak Djokovic")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Value2 -eq "Rolland Garros") {
            $cell.Value2 = "Roland Garros"
        }
    }
}

# --- Update the saved selection / active cell on each sheet, and make
#     "Roger Federer" the active tab again (it was "Novak Djokovic"). ---
# Select the non-active sheets first, then select on the sheet that should
# end up being the active tab last, so it "wins" the tab-selected flag.
$wsNadal = $wb.Worksheets.Item("Rafael Nadal")
$wsNadal.Range("B2").Select()

$wsDjokovic = $wb.Worksheets.Item("Novak Djokovic")
$wsDjokovic.Range("B21").Select()

$wsFederer = $wb.Worksheets.Item("Roger Federer")
$wsFederer.Range("L9").Select()
